$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new article row (row 7): Address in column A, Description in column B
$ws.Range("A7").Value = "article-6.html"
$ws.Range("B7").Value = "How Are WIS students coping with this new Pandemic?"

# Widen column B to fit the new, longer description text
$ws.Columns.Item(2).ColumnWidth = 48.44140625
